$wb = $excel.ActiveWorkbook
$st = $wb.Styles.Add("MyStyle")
$st.Font.Name = "Wingdings"
Write-Output "added"
